# Auto-generated edit script: wraps specific CV placeholder-candidate paragraphs
# into Jinja-style {{ field }} template tags, matching the target diff.
# Each paragraph is located by its unique literal text via Find, then the whole
# paragraph's Range is replaced with an equivalent WordprocessingML fragment
# (same pPr/paraId/rsid attributes) whose runs spell out '{{', the field name
# (wrapped in spellStart/spellEnd proofErr marks, as Word does for non-dictionary
# words), and '}}'.

$d = $word.ActiveDocument

# --- Paragraph 1: '3.0' ---
$rng = $d.Content
$found = $rng.Find.Execute('3.0', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "paragraph not found: 3.0" }
$para = $rng.Paragraphs(1)
$xml0 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="4C9ACE92" w14:textId="737D4B82" w:rsidR="007446CD" w:rsidRDefault="00B360D2"><w:pPr><w:pStyle w:val="ListeMaddemi"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>gpa</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$para.Range.InsertXML($xml0)

# --- Paragraph 2: 'Car washing place 3 years, washer machine, 8 months 2001-2004' ---
$rng = $d.Content
$found = $rng.Find.Execute('Car washing place 3 years, washer machine, 8 months 2001-2004', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "paragraph not found: Car washing place 3 years, washer machine, 8 months 2001-2004" }
$para = $rng.Paragraphs(1)
$xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="6B1174E2" w14:textId="77777777" w:rsidR="00B360D2" w:rsidRPr="00B360D2" w:rsidRDefault="00B360D2" w:rsidP="00B360D2"><w:pPr><w:pStyle w:val="Blm"/><w:spacing w:after="0"/><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>work_experience</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>}}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$para.Range.InsertXML($xml1)

# --- Paragraph 3: 'Nvidia, founder, 21 years 3 months' ---
$rng = $d.Content
$found = $rng.Find.Execute('Nvidia, founder, 21 years 3 months', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "paragraph not found: Nvidia, founder, 21 years 3 months" }
$para = $rng.Paragraphs(1)
$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="21B7545C" w14:textId="77777777" w:rsidR="00B360D2" w:rsidRPr="00B360D2" w:rsidRDefault="00B360D2" w:rsidP="00B360D2"><w:pPr><w:pStyle w:val="Blm"/><w:spacing w:after="0"/><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>work_experience</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>}}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$para.Range.InsertXML($xml2)

# --- Paragraph 4: 'can setup arch' ---
$rng = $d.Content
$found = $rng.Find.Execute('can setup arch', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "paragraph not found: can setup arch" }
$para = $rng.Paragraphs(1)
$xml3 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="4A077FC3" w14:textId="76D78306" w:rsidR="007446CD" w:rsidRDefault="00B3313A" w:rsidP="000108FC"><w:pPr><w:pStyle w:val="ListeMaddemi"/></w:pPr><w:r><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>qualifications</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$para.Range.InsertXML($xml3)

# --- Paragraph 5: '5+ years in spring systems' ---
$rng = $d.Content
$found = $rng.Find.Execute('5+ years in spring systems', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "paragraph not found: 5+ years in spring systems" }
$para = $rng.Paragraphs(1)
$xml4 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="48E716AE" w14:textId="77777777" w:rsidR="000108FC" w:rsidRDefault="000108FC"><w:pPr><w:pStyle w:val="ListeMaddemi"/></w:pPr><w:r><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>qualifications</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$para.Range.InsertXML($xml4)

# --- Paragraph 6: 'can spell phyton correctly' ---
$rng = $d.Content
$found = $rng.Find.Execute('can spell phyton correctly', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "paragraph not found: can spell phyton correctly" }
$para = $rng.Paragraphs(1)
$xml5 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="7A476C18" w14:textId="77777777" w:rsidR="000108FC" w:rsidRDefault="000108FC"><w:pPr><w:pStyle w:val="ListeMaddemi"/></w:pPr><w:r><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>qualifications</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$para.Range.InsertXML($xml5)

# --- Paragraph 7: 'gateway management expert, stone to stone' ---
$rng = $d.Content
$found = $rng.Find.Execute('gateway management expert, stone to stone', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "paragraph not found: gateway management expert, stone to stone" }
$para = $rng.Paragraphs(1)
$xml6 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="6E66F93B" w14:textId="5328CCB8" w:rsidR="00B3313A" w:rsidRDefault="00B3313A"><w:pPr><w:pStyle w:val="ListeMaddemi"/></w:pPr><w:r><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>qualifications</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$para.Range.InsertXML($xml6)
